$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated likelihood-table values (rows 2-16) after re-running the model:
# D = R (raw measurement), E = C (raw measurement), F = R|L_mean, G = R|L_var,
# H = C|L_mean, I = C|L_var, J = P(R|L), K = P(C|L), L = P(L) * P(R|L) * P(C|L)
$newValues = @{
    "2" = @{ "D" = 1.027; "E" = 5.3; "F" = 1.23625; "G" = 0.09761414583333339; "H" = 3.183333333333334; "I" = 1.338055555555555; "J" = 1.020353483153802; "K" = 0.06465291590195688; "L" = 0.03298441396830576 }
    "3" = @{ "D" = 1.239; "E" = 4; "F" = 1.23625; "G" = 0.09761414583333339; "H" = 3.183333333333334; "I" = 1.338055555555555; "J" = 1.27684112895035; "K" = 0.2688050773771029; "L" = 0.1716106892328831 }
    "4" = @{ "D" = 1.166; "E" = 2.5; "F" = 1.23625; "G" = 0.09761414583333339; "H" = 3.183333333333334; "I" = 1.338055555555555; "J" = 1.245017365128281; "K" = 0.2896640519617951; "L" = 0.1803183873729277 }
    "5" = @{ "D" = 0.9089999999999999; "E" = 3; "F" = 1.23625; "G" = 0.09761414583333339; "H" = 3.183333333333334; "I" = 1.338055555555555; "J" = 0.7377704624374438; "K" = 0.3405794120351475; "L" = 0.1256347151569218 }
    "6" = @{ "D" = 1.1865; "E" = 2.5; "F" = 1.23625; "G" = 0.09761414583333339; "H" = 3.183333333333334; "I" = 1.338055555555555; "J" = 1.260804627539021; "K" = 0.2896640519617951; "L" = 0.1826048885725673 }
    "7" = @{ "D" = 1.89; "E" = 1.8; "F" = 1.23625; "G" = 0.09761414583333339; "H" = 3.183333333333334; "I" = 1.338055555555555; "J" = 0.143023266457987; "K" = 0.1687024282500931; "L" = 0.01206418617386125 }
    "8" = @{ "D" = 2.171; "E" = 5; "F" = 1.694111111111111; "G" = 0.2561913765432099; "H" = 3.638888888888889; "I" = 1.517098765432099; "J" = 0.5056669539171873; "K" = 0.1758859568782206; "L" = 0.04446985802570978 }
    "9" = @{ "D" = 1.8415; "E" = 2.9; "F" = 1.694111111111111; "G" = 0.2561913765432099; "H" = 3.638888888888889; "I" = 1.517098765432099; "J" = 0.7554661938984805; "K" = 0.2705568796570144; "L" = 0.1021982880537669 }
    "10" = @{ "D" = 1.031; "E" = 2.95; "F" = 1.694111111111111; "G" = 0.2561913765432099; "H" = 3.638888888888889; "I" = 1.517098765432099; "J" = 0.3341374974572674; "K" = 0.2769980447648336; "L" = 0.04627771673913881 }
    "11" = @{ "D" = 1.742; "E" = 3; "F" = 1.694111111111111; "G" = 0.2561913765432099; "H" = 3.638888888888889; "I" = 1.517098765432099; "J" = 0.7846644510953596; "K" = 0.2831256128309231; "L" = 0.1110793017915068 }
    "12" = @{ "D" = 2.302; "E" = 1.7; "F" = 1.694111111111111; "G" = 0.2561913765432099; "H" = 3.638888888888889; "I" = 1.517098765432099; "J" = 0.3831915621817934; "K" = 0.0938261323161729; "L" = 0.01797669110785497 }
    "13" = @{ "D" = 1.2415; "E" = 4.9; "F" = 1.694111111111111; "G" = 0.2561913765432099; "H" = 3.638888888888889; "I" = 1.517098765432099; "J" = 0.5284350661159578; "K" = 0.1917625863465881; "L" = 0.05066703749731317 }
    "14" = @{ "D" = 1.07; "E" = 4.25; "F" = 1.694111111111111; "G" = 0.2561913765432099; "H" = 3.638888888888889; "I" = 1.517098765432099; "J" = 0.3685327064981622; "K" = 0.2863842059988284; "L" = 0.05277097326753772 }
    "15" = @{ "D" = 2.446; "E" = 2.55; "F" = 1.694111111111111; "G" = 0.2561913765432099; "H" = 3.638888888888889; "I" = 1.517098765432099; "J" = 0.2614866445343001; "K" = 0.2191255786365801; "L" = 0.02864920614465813 }
    "16" = @{ "D" = 1.402; "E" = 5.5; "F" = 1.694111111111111; "G" = 0.2561913765432099; "H" = 3.638888888888889; "I" = 1.517098765432099; "J" = 0.6672724976589723; "K" = 0.1034255012828138; "L" = 0.0345064962813072 }
}

foreach ($rowKey in $newValues.Keys) {
    $rowVals = $newValues[$rowKey]
    foreach ($colKey in $rowVals.Keys) {
        $addr = "$colKey$rowKey"
        $ws.Range($addr).Value = $rowVals[$colKey]
    }
}
